$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactUspage")

# --- Row 1: headers (left to right so shared-string order matches) ---
$ws.Range("C1").Value = "firstName"
$ws.Range("D1").Value = "LastName"
$ws.Range("E1").Value = "email"
$ws.Range("F1").Value = "phone"
$ws.Range("G1").Value = "message"

# --- Row 2: first contact-us submission ---
$ws.Range("C2").Value = "Shilpa"
$ws.Range("D2").Value = "N"
$ws.Range("E2").Value = "nshilpamurthy@gmail.com"
$ws.Range("F2").Value = 7483067146
$ws.Range("G2").Value = "hi hello "

# --- Row 3: repeat submission (reuses the same shared strings) ---
$ws.Range("C3").Value = "Shilpa"
$ws.Range("E3").Value = "nshilpamurthy@gmail.com"
$ws.Range("F3").Value = 7483067146
$ws.Range("G3").Value = "hi hello "

# --- Hyperlinks for the email cells (E3 first, then E2, to mirror rId order) ---
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:nshilpamurthy@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:nshilpamurthy@gmail.com")

# --- Column sizing to roughly match the new data layout (closest the
#     engine's column-width quantization allows to the recorded widths) ---
$ws.Columns.Item(1).ColumnWidth = 19.666666667
$ws.Columns.Item(2).ColumnWidth = 80
$ws.Columns.Item(3).ColumnWidth = 9
$ws.Columns.Item(4).ColumnWidth = 8.833333333
$ws.Columns.Item(5).ColumnWidth = 24.5
$ws.Columns.Item(6).ColumnWidth = 10.166666667
$ws.Columns.Item(7).ColumnWidth = 7.833333333

$ws.Range("B18").Select() | Out-Null
